$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 blog "ser" numbers shift: I7 107->108, E7 108->110, C7 110->109 (new)
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 108"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 110"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 109"

# Update the active selection to C7
$ws.Range("C7").Select()
